$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width: 28 -> 7 (XML width accounts for a fixed 0.8333 padding
# offset between the COM ColumnWidth property and the stored XML width)
$ws.Range("E1").ColumnWidth = 6.1666666666666667

# Name corrections to uppercase
$ws.Range("A4").Value = "PEDRO CUEVAS"
$ws.Range("A22").Value = "PEDRO CUEVAS"
$ws.Range("A30").Value = "PEDRO CUEVAS"
$ws.Range("A48").Value = "PEDRO CUEVAS"

# D4 / D30 cedula number (previously empty)
$ws.Range("D4").Value = "0-472963882991"
$ws.Range("D30").Value = "0-472963882991"

# Monto corrections 163.17 -> 166.87
$ws.Range("G4").Value = 166.87
$ws.Range("G8").Value = 166.87
$ws.Range("G9").Value = 166.87
$ws.Range("G12").Value = 166.87
$ws.Range("C13").Value = 166.87
$ws.Range("C17").Value = 166.87
$ws.Range("B20").Value = 166.87
$ws.Range("G20").Value = 166.87
$ws.Range("G30").Value = 166.87
$ws.Range("G34").Value = 166.87
$ws.Range("G35").Value = 166.87
$ws.Range("G38").Value = 166.87
$ws.Range("C39").Value = 166.87
$ws.Range("C43").Value = 166.87
$ws.Range("B46").Value = 166.87
$ws.Range("G46").Value = 166.87

# Salario por hora / diario corrections
$ws.Range("C8").Value = 4.51
$ws.Range("C34").Value = 4.51
$ws.Range("C9").Value = 36.08
$ws.Range("C35").Value = 36.08

# Banco General -> BANCO GENERAL (uppercase)
$ws.Range("D15").Value = "BANCO GENERAL"
$ws.Range("C20").Value = "BANCO GENERAL"
$ws.Range("D41").Value = "BANCO GENERAL"
$ws.Range("C46").Value = "BANCO GENERAL"
